$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newText = "Upon meeting the commander of the Chinese team for the first time, he looks slightly discouraged when he sees your rank. You explain that you have taken the place of your commander because he had unexpected business come up. He has proposed to set up two aid stations to help victims of a natural distaster, one American and one Chinese. You think it would be better to have one large American aid station, so you ask what benefit his plan has."

# Insert a new row at row 14, shifting rows 14-17 down to 15-18
$ws.Rows.Item(14).Insert()

# Set the new row's cells
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = $newText

$ws.Range("A15").Value = 14

$ws.Range("A16").Value = 15

$ws.Range("A17").Value = 16

$ws.Range("A18").Value = 17

# Update sheet view state (scroll position + selection)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B14").Select()
